$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Delete row 5 ("ford tourneo custom 2024" duplicate). Everything below shifts up by one.
$ws.Rows.Item(5).Delete()

# After the shift, the former row 17 (duplicate "Honda CR V 2024") is now row 16; remove it too.
$ws.Rows.Item(16).Delete()
